# Scheduled market-price refresh: update H/I/J/K/L/M/N price & profit
# columns across the per-job Leve Profit sheets (values sourced from the
# latest Universalis market snapshot).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 42571.43
$ws.Range("J3").Value = 42571.43
$ws.Range("L3").Value = 42571.43
$ws.Range("N3").Value = -42799.43
$ws.Range("H21").Value = 57634.5
$ws.Range("J21").Value = 42500
$ws.Range("L21").Value = 42500
$ws.Range("N21").Value = -43436
$ws.Range("H23").Value = 57634.5
$ws.Range("J23").Value = 42500
$ws.Range("L23").Value = 42500
$ws.Range("N23").Value = -42968
$ws.Range("H102").Value = 42571.43
$ws.Range("J102").Value = 42571.43
$ws.Range("L102").Value = 42571.43
$ws.Range("N102").Value = -49061.43
$ws.Range("H116").Value = 5754.5557
$ws.Range("I116").Value = 6284.4287
$ws.Range("J116").Value = 3900
$ws.Range("K116").Value = 6284.4287
$ws.Range("L116").Value = 3900
$ws.Range("M116").Value = -2842.4287
$ws.Range("N116").Value = -10784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1068.7931
$ws.Range("I2").Value = 708.86365
$ws.Range("J2").Value = 2200
$ws.Range("K2").Value = 708.86365
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = -595.86365
$ws.Range("N2").Value = -2426
$ws.Range("H32").Value = 4730.3784
$ws.Range("I32").Value = 3347.9285
$ws.Range("K32").Value = 3347.9285
$ws.Range("M32").Value = -3060.9285
$ws.Range("H45").Value = 4524.115
$ws.Range("I45").Value = 5390.4287
$ws.Range("K45").Value = 5390.4287
$ws.Range("M45").Value = -5013.4287
$ws.Range("H116").Value = 1068.7931
$ws.Range("I116").Value = 708.86365
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 708.86365
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = 1585.13635
$ws.Range("N116").Value = -6788
$ws.Range("H132").Value = 1975.5146
$ws.Range("I132").Value = 1119.1837
$ws.Range("J132").Value = 4183.9473
$ws.Range("K132").Value = 3357.5511
$ws.Range("L132").Value = 12551.8419
$ws.Range("M132").Value = -827.5511000000001
$ws.Range("N132").Value = -17611.8419

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1068.7931
$ws.Range("I3").Value = 708.86365
$ws.Range("J3").Value = 2200
$ws.Range("K3").Value = 708.86365
$ws.Range("L3").Value = 2200
$ws.Range("M3").Value = -594.86365
$ws.Range("N3").Value = -2428
$ws.Range("H22").Value = 827.17645
$ws.Range("I22").Value = 783.5
$ws.Range("J22").Value = 851
$ws.Range("K22").Value = 783.5
$ws.Range("L22").Value = 851
$ws.Range("M22").Value = -610.5
$ws.Range("N22").Value = -1197
$ws.Range("H94").Value = 1607.0714
$ws.Range("I94").Value = 1299.4
$ws.Range("J94").Value = 2376.25
$ws.Range("K94").Value = 1299.4
$ws.Range("L94").Value = 2376.25
$ws.Range("M94").Value = -848.4000000000001
$ws.Range("N94").Value = -3278.25
$ws.Range("H100").Value = 29800
$ws.Range("J100").Value = 29800
$ws.Range("L100").Value = 29800
$ws.Range("N100").Value = -31964
$ws.Range("H130").Value = 49780
$ws.Range("J130").Value = 49780
$ws.Range("L130").Value = 49780
$ws.Range("N130").Value = -59820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6174483.5
$ws.Range("I58").Value = 12821624
$ws.Range("J58").Value = 2138.2856
$ws.Range("K58").Value = 12821624
$ws.Range("L58").Value = 2138.2856
$ws.Range("M58").Value = -12821421
$ws.Range("N58").Value = -2544.2856
$ws.Range("H99").Value = 4468033.5
$ws.Range("I99").Value = 2930.1904
$ws.Range("J99").Value = 17863342
$ws.Range("K99").Value = 2930.1904
$ws.Range("L99").Value = 17863342
$ws.Range("M99").Value = -1432.1904
$ws.Range("N99").Value = -17866338
$ws.Range("H122").Value = 1253.4
$ws.Range("I122").Value = 1187.2
$ws.Range("K122").Value = 3561.6
$ws.Range("M122").Value = -1111.6
$ws.Range("H126").Value = 4468033.5
$ws.Range("I126").Value = 2930.1904
$ws.Range("J126").Value = 17863342
$ws.Range("K126").Value = 8790.5712
$ws.Range("L126").Value = 53590026
$ws.Range("M126").Value = -6320.5712
$ws.Range("N126").Value = -53594966
$ws.Range("H136").Value = 6174483.5
$ws.Range("I136").Value = 12821624
$ws.Range("J136").Value = 2138.2856
$ws.Range("K136").Value = 38464872
$ws.Range("L136").Value = 6414.8568
$ws.Range("M136").Value = -38462322
$ws.Range("N136").Value = -11514.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 983.3333
$ws.Range("J92").Value = 1233.3334
$ws.Range("L92").Value = 3700.0002
$ws.Range("N92").Value = -6196.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9983.75
$ws.Range("J57").Value = 9983.75
$ws.Range("L57").Value = 9983.75
$ws.Range("N57").Value = -11623.75
$ws.Range("H122").Value = 4594993
$ws.Range("I122").Value = 4987462.5
$ws.Range("J122").Value = 4169817.2
$ws.Range("K122").Value = 14962387.5
$ws.Range("L122").Value = 12509451.6
$ws.Range("M122").Value = -14959937.5
$ws.Range("N122").Value = -12514351.6
$ws.Range("H126").Value = 5348.8066
$ws.Range("I126").Value = 8247.467000000001
$ws.Range("J126").Value = 2631.3125
$ws.Range("K126").Value = 24742.401
$ws.Range("L126").Value = 7893.9375
$ws.Range("M126").Value = -22272.401
$ws.Range("N126").Value = -12833.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 502442.22
$ws.Range("I82").Value = 715577.0600000001
$ws.Range("J82").Value = 129456.25
$ws.Range("K82").Value = 715577.0600000001
$ws.Range("L82").Value = 129456.25
$ws.Range("M82").Value = -715216.0600000001
$ws.Range("N82").Value = -130178.25
$ws.Range("H85").Value = 502442.22
$ws.Range("I85").Value = 715577.0600000001
$ws.Range("J85").Value = 129456.25
$ws.Range("K85").Value = 715577.0600000001
$ws.Range("L85").Value = 129456.25
$ws.Range("M85").Value = -714329.0600000001
$ws.Range("N85").Value = -131952.25
$ws.Range("H132").Value = 11757017
$ws.Range("I132").Value = 15283141
$ws.Range("J132").Value = 3270.7778
$ws.Range("K132").Value = 45849423
$ws.Range("L132").Value = 9812.3334
$ws.Range("M132").Value = -45846893
$ws.Range("N132").Value = -14872.3334
$ws.Range("H136").Value = 7893.2163
$ws.Range("I136").Value = 8491.723
$ws.Range("J136").Value = 7326.2104
$ws.Range("K136").Value = 25475.169
$ws.Range("L136").Value = 21978.6312
$ws.Range("M136").Value = -22925.169
$ws.Range("N136").Value = -27078.6312

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 22225466
$ws.Range("I81").Value = 3000.5
$ws.Range("J81").Value = 28574742
$ws.Range("K81").Value = 6001
$ws.Range("L81").Value = 57149484
$ws.Range("M81").Value = -4940
$ws.Range("N81").Value = -57151606
$ws.Range("H84").Value = 22225466
$ws.Range("I84").Value = 3000.5
$ws.Range("J84").Value = 28574742
$ws.Range("K84").Value = 30005
$ws.Range("L84").Value = 285747420
$ws.Range("M84").Value = -24701
$ws.Range("N84").Value = -285758028
$ws.Range("H113").Value = 1080.9231
$ws.Range("J113").Value = 918.0714
$ws.Range("L113").Value = 2754.2142
$ws.Range("N113").Value = -7094.2142
